$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.341.96'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.926.77'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.44'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.87'
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.99'
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.45'
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").Value = '3.412.39'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '61.338.55'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '2.927.94'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.68'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '431.55'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("E21").Value = '  -1.27%  '
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.77'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.87'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("E25").Value = '  -1.96%  '
$ws.Range("E26").Value = '  -2.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -5.11%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.109'
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.57'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '0.0₃0883'
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.97'
$ws.Range("E37").Value = '  -2.20%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.123'
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '42.12'
$ws.Range("E41").Value = '  +5.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.280'
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").Value = '2.698.49'
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '133.89'
$ws.Range("E45").Value = '  +2.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '364.24'
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.59'
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.124'
$ws.Range("E51").Value = '  -2.15%  '

# Restore default (no explicit number format) for cells forced to text,
# matching the original workbook where these cells had no style override.
foreach ($addr in @("D5", "D6", "D9", "D13", "D18", "D19", "D23", "D24", "D28", "D31", "D32", "D37", "D39", "D41", "D42", "D45", "D46", "D48", "D51")) {
    $ws.Range($addr).ClearFormats()
}
